# Finish the report creator tool: add the extra rows gathered by the
# tool, normalize the formatting of the whole table body and widen the
# name column so the longer labels are readable.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new rows to the report.
$ws.Range("A10").Value = "Chicken gun rooster rudy"
$ws.Range("B10").Value = 32
$ws.Range("C10").Value = 3257

$ws.Range("A11").Value = "Tio chico"
$ws.Range("B11").Value = 328
$ws.Range("C11").Value = 312358

$ws.Range("A12").Value = "Total"
$ws.Range("B12").Value = "fex"
$ws.Range("C12").Value = "fex"

$ws.Range("A13").Value = "pitsecato"
$ws.Range("B13").Value = "fex"
$ws.Range("C13").Value = "fex"

$ws.Range("A14").Value = "fortfire leaks"
$ws.Range("B14").Value = "fex"
$ws.Range("C14").Value = "fex"

# Make B:C match column A's (Normal) cell style for the whole body, old
# rows and new ones alike, in one pass so they all pick up the same style.
$ws.Range("B5:C14").Style = "Normal"

# Widen column A and move the active selection to C13 like the source file.
$ws.Columns("A").ColumnWidth = 23
$ws.Range("C13").Select()
